$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the "actual value" (F) and "execution result" (G) columns for all
# data rows - these held redundant/duplicate data that's being cleaned up.
$ws.Range("F2:G20").ClearContents()

# Update the active selection to reflect where the user ended up (F19).
$ws.Range("F19").Select()
